$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01213623907085998
$ws.Range("C2").Value = 0.0129912521216869
$ws.Range("D2").Value = 0.01256374559627344
$ws.Range("E2").Value = 0.0004275065254134638

$ws.Range("B3").Value = 0.2995169082125604
$ws.Range("C3").Value = 0.3220064724919094
$ws.Range("D3").Value = 0.3107616903522349
$ws.Range("E3").Value = 0.0112447821396745

$ws.Range("B4").Value = 0.02332727158713237
$ws.Range("C4").Value = 0.02497489959839358
$ws.Range("D4").Value = 0.02415108559276297
$ws.Range("E4").Value = 0.000823814005630601
